# Apply "latest code commited on 15th sept" changes:
#  - Update a handful of INC/WO reference values (shared strings get
#    renumbered/pruned automatically by the engine when old values are
#    no longer referenced).
#  - Move the active sheet / selection cursor on several sheets to match
#    the new view state (Create_Inc, Update_Inc, WO_creation, WO_search).

$wb = $excel.ActiveWorkbook

# --- Cell value updates -----------------------------------------------

$wsSearchInc = $wb.Worksheets.Item("Search_Inc")
$wsSearchInc.Range("A2").Value = "INC000000533091"

$wsCreateInc = $wb.Worksheets.Item("Create_Inc")
$wsCreateInc.Range("W2").Value = "INC000000533091"

$wsUpdateInc = $wb.Worksheets.Item("Update_Inc")
$wsUpdateInc.Range("A2").Value = "INC000000533091"
$wsUpdateInc.Range("A3").Value = "INC000000533091"
$wsUpdateInc.Range("A4").Value = "INC000000533091"

$wsWoSearch = $wb.Worksheets.Item("WO_search")
$wsWoSearch.Range("A2").Value = "WO0000000160163 "

# --- View / selection state --------------------------------------------
# Order matters: the last sheet selected becomes the active tab, which
# must end up being WO_search (activeTab index 5).

$wsCreateInc.Range("D17").Select()

$wsUpdateInc.Range("A3:A4").Select()

$wsWoCreation = $wb.Worksheets.Item("WO_creation")
$wsWoCreation.Range("H8").Select()

$wsWoSearch.Range("A6").Select()
